# Add four new attendance-tracking columns (lates, absents, latetoday,
# absenttoday) to the 12STEM roster sheet.
#
# Column layout after the edit:
#   A-F : existing roster columns (unchanged)
#   G   : lates        (number, default 0)
#   H   : absents      (number, default 0)
#   I   : latetoday    (boolean, default FALSE)
#   J   : absenttoday  (boolean, default TRUE)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 23

# --- Header row -------------------------------------------------------
$ws.Range("G1").Value = "lates"
$ws.Range("H1").Value = "absents"
$ws.Range("I1").Value = "latetoday"
$ws.Range("J1").Value = "absenttoday"

# G1:I1 wrap their header text (matches the new wrapped/left-aligned style);
# J1 keeps the plain left-aligned header style already used by columns A-F.
$ws.Range("G1:I1").WrapText = $true
$ws.Range("G1:J1").HorizontalAlignment = -4131

$ws.Rows.Item(1).RowHeight = 30.75

# --- Data rows ----------------------------------------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = 0        # G: lates
    $ws.Cells.Item($r, 8).Value = 0        # H: absents
    $ws.Cells.Item($r, 9).Value = $false   # I: latetoday
    $ws.Cells.Item($r, 10).Value = $true   # J: absenttoday
}

# --- Selection ------------------------------------------------------------
$null = $ws.Range("G2:J23").Select()
